# Update CompResult / Sheet3 with refreshed computation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Updated K-svmeans (column E) figures for the first results block (rows 3-19).
$ws.Range("E3").Value = 71
$ws.Range("E8").Value = 82
$ws.Range("E9").Value = 85.5
$ws.Range("E12").Value = 95.5
$ws.Range("E18").Value = 67

# Move the active selection to reflect where the user left off reviewing the update.
$ws.Range("E9").Select()
